$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Лист1 -> Sheet1)
$ws.Name = "Sheet1"

# Clear any previous content in the used range
$ws.Cells.Clear()

# Write the new user-list data (row 2 and row 3; row 1 stays empty, matching the dimension A2:G3)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "admin"
$ws.Range("C2").Value = "password"
$ws.Range("D2").Value = "Shakira Regalado"
$ws.Range("E2").Value = "shakiraregalado@gmail.com"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "Fluffy"

$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "kira"
$ws.Range("C3").Value = "pass"
$ws.Range("D3").Value = "Shakira"
$ws.Range("E3").Value = "shakira@gmail.com"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = "Pink"
